# Fruta / hortaliza, semanal
# Re-orders the weekly price records held in rows 3-9 (columns D, M, N, O, P, Q, S, T)
# so that the data lines up with the correct dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row, after the re-shuffle described in the commit.
$rows = @{
    3 = @{ D = 44418; M = 240; N = 10000; O = 11000; P = 10500; Q = "`$/bandeja 10 kilos"; S = 1050; T = 10 }
    4 = @{ D = 44489; M = 300; N = 26000; O = 27000; P = 26500; Q = "`$/bandeja 18 kilos"; S = 1472; T = 18 }
    5 = @{ D = 44291; M = 200; N = 17000; O = 18000; P = 17500; Q = "`$/bandeja 18 kilos"; S = 972;  T = 18 }
    6 = @{ D = 44487; M = 300; N = 14000; O = 15000; P = 14500; Q = "`$/bandeja 10 kilos"; S = 1450; T = 10 }
    7 = @{ D = 44323; M = 270; N = 21000; O = 22000; P = 21500; Q = "`$/bandeja 18 kilos"; S = 1194; T = 18 }
    8 = @{ D = 44491; M = 300; N = 14000; O = 15000; P = 14500; Q = "`$/bandeja 10 kilos"; S = 1450; T = 10 }
    9 = @{ D = 44307; M = 250; N = 19000; O = 20000; P = 19500; Q = "`$/bandeja 18 kilos"; S = 1083; T = 18 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
